$wb = $excel.ActiveWorkbook

# Rename Sheet1 to ReleaseDates
$wb.Worksheets.Item("Sheet1").Name = "ReleaseDates"

# Delete Sheet2 and Sheet3
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()
$excel.DisplayAlerts = $true
